$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Add the new "AchievementType" lookup-table sheet at the very end
#    of the workbook (after "PartnerPermission").
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$achievementType = $wb.Worksheets.Add($null, $lastSheet)
$achievementType.Name = "AchievementType"

$achievementType.Columns.Item(2).ColumnWidth = 23.28515625

$achievementType.Range("A1").Value = "AchievementType"
$achievementType.Range("A2").Value = "Id"
$achievementType.Range("B2").Value = "Name"
$achievementType.Range("C2").Formula = "=CONCATENATE(`"delete from `",`$A`$1,`"; dbcc checkident (`",`$A`$1,`", reseed, 0);`")"

$achievementType.Range("A3").Value = 1
$achievementType.Range("B3").Value = "Transakcija"
$achievementType.Range("C3").Formula = "=CONCATENATE(`"insert into `",`$A`$1,`"(`",`$B`$2,`") values(N'`",B3,`"');`")"

# ---------------------------------------------------------------------
# 2) Fix up the existing "Genders" sheet: translate Male/Female into
#    Serbian (Muski/Zenski) and fix the off-by-one formula bug that
#    referenced the row above instead of the current row.
# ---------------------------------------------------------------------
$genders = $wb.Worksheets.Item("Genders")

$genders.Range("B3").Value = "Muški"
$genders.Range("C3").Formula = "=CONCATENATE(`"insert into `",`$A`$1,`"(`",`$B`$2,`") values(N'`",B3,`"');`")"

$genders.Range("B4").Value = "Ženski"
$genders.Range("C4").Formula = "=CONCATENATE(`"insert into `",`$A`$1,`"(`",`$B`$2,`") values(N'`",B4,`"');`")"

$genders.Activate()
$genders.Range("D10").Select()

# ---------------------------------------------------------------------
# 3) Back to "AchievementType": finish populating the remaining rows.
# ---------------------------------------------------------------------
$achievementType.Range("A4").Value = 2
$achievementType.Range("B4").Value = "Popunjavanje prvi put"
$achievementType.Range("C4").Formula = "=CONCATENATE(`"insert into `",`$A`$1,`"(`",`$B`$2,`") values(N'`",B4,`"');`")"

$achievementType.Range("A5").Value = 3
$achievementType.Range("B5").Value = "Manuelno"
$achievementType.Range("C5").Formula = "=CONCATENATE(`"insert into `",`$A`$1,`"(`",`$B`$2,`") values(N'`",B5,`"');`")"

$achievementType.Activate()
$achievementType.Range("C9").Select()
